# Adds a "Meta description" paragraph right after the document title, and
# replaces the old duplicate title/description block at the end of the
# document with an image-generation prompt paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new paragraph right after the title ("Play Dancing in Rio
#    Slot for Free - Review") containing a bold "Meta description" run
#    followed by a plain run with the description text.
#    We build the new paragraph (plus a throw-away trailing paragraph
#    mark, required so InsertXML treats it as a standalone paragraph
#    break rather than merging into the following heading) via
#    Range.InsertXML, then delete the extra empty paragraph it leaves
#    behind.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$insertPos = $titlePara.Range.End
$insertTarget = $d.Range($insertPos, $insertPos)

$metaXml = '<?xml version="1.0" standalone="yes"?>' + `
  '<?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:r/>' + `
              '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
              '<w:r><w:t>: Read our review of Dancing in Rio to discover its features and see if you want to play it for free. Includes a progressive jackpot and up to 50 free spins.</w:t></w:r>' + `
            '</w:p>' + `
            '<w:p><w:pPr><w:rPr></w:rPr></w:pPr></w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$insertTarget.InsertXML($metaXml)

# Remove the extra empty paragraph left behind by the fragment's trailing
# paragraph mark (it is now paragraph #3, right before "Gameplay Features").
$extraPara = $d.Paragraphs.Item(3)
if ($extraPara.Range.Text -eq "" -or $extraPara.Range.Text -eq "`r") {
    $extraPara.Range.Delete()
}

# ---------------------------------------------------------------------
# 2) Remove the duplicate bold title paragraph ("Play Dancing in Rio
#    Slot for Free - Review") that used to sit near the end of the
#    document, just before the italic description paragraph.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
if ($dupTitlePara.Range.Text -match "Play Dancing in Rio Slot for Free - Review") {
    $dupTitlePara.Range.Delete()
}

# ---------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    image-generation prompt, keeping its leading empty run and italic
#    formatting intact.
# ---------------------------------------------------------------------
$oldText = "Read our review of Dancing in Rio to discover its features and see if you want to play it for free. Includes a progressive jackpot and up to 50 free spins."
$newText = "Create a feature image for the Dancing in Rio game in a cartoon style. The image should feature a happy Maya warrior with glasses. The Maya warrior should be surrounded by other happy carnival dancers in brightly colored costumes, all dancing in the streets with confetti raining down on them. The image should capture the joyful and vibrant vibes of the Rio Carnival. The background should be the iconic Christ the Redeemer statue, adding a touch of authenticity to the image. The overall feel of the image should be fun and lively, inviting players to join in on the carnival festivities and try their luck at this exciting slot game."

$finalCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($finalCount)
$lastPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                              $true, 1, $false, $newText, 2)

Write-Output "Edit complete."
